$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) text updates ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) value updates for rows 2-35 ---
$ws.Range("C2").Value = 4547.50930098406
$ws.Range("C3").Value = 3972.630273980753
$ws.Range("C4").Value = 17288.8595992193
$ws.Range("C5").Value = 1280.225469721551
$ws.Range("C6").Value = 4633.590358399045
$ws.Range("C7").Value = 4355.934938677345
$ws.Range("C8").Value = 2217.474008566157
$ws.Range("C9").Value = 17610.30663334184
$ws.Range("C10").Value = 1263.452411343738
$ws.Range("C11").Value = 2024.117324382548
$ws.Range("C12").Value = 4479.398934239905
$ws.Range("C13").Value = 2264.394087033834
$ws.Range("C14").Value = 2094.024217383061
$ws.Range("C15").Value = 4394.543881413723
$ws.Range("C16").Value = 18254.09644617555
$ws.Range("C17").Value = 1291.622214254295
$ws.Range("C18").Value = 2379.668184479739
$ws.Range("C19").Value = 2201.396847776877
$ws.Range("C20").Value = 4699.493713911862
$ws.Range("C21").Value = 16764.42871195103
$ws.Range("C22").Value = 1291.415042301529
$ws.Range("C23").Value = 2497.68592515536
$ws.Range("C24").Value = 2612.856880840196
$ws.Range("C25").Value = 16146.07242861928
$ws.Range("C26").Value = 711.3043470146426
$ws.Range("C27").Value = 4861.287098802361
$ws.Range("C28").Value = 2735.187532014817
$ws.Range("C29").Value = 14093.81249338665
$ws.Range("C30").Value = 731.9993357350996
$ws.Range("C31").Value = 4944.191641077407
$ws.Range("C32").Value = 5176.058803160127
$ws.Range("C33").Value = 2886.897484630703
$ws.Range("C34").Value = 809.9545825255682
$ws.Range("C35").Value = 5089.61202008711
